$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 138, shifting rows 138:157 down to 139:158
$ws.Rows.Item(138).Insert()

# Populate the new row 138 with the new week's data
$ws.Cells.Item(138, 1).Value = 9
$ws.Cells.Item(138, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(138, 3).Value = "Metropolitana"
$ws.Cells.Item(138, 4).Value = 44474
$ws.Cells.Item(138, 5).Value = 13
$ws.Cells.Item(138, 6).Value = 300000001
$ws.Cells.Item(138, 7).Value = "Rabanito"
$ws.Cells.Item(138, 8).Value = "Sin especificar"
$ws.Cells.Item(138, 9).Value = "Primera"
$ws.Cells.Item(138, 10).Value = 7900
$ws.Cells.Item(138, 11).Value = 3500
$ws.Cells.Item(138, 12).Value = 4000
$ws.Cells.Item(138, 13).Value = 3747
$ws.Cells.Item(138, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(138, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(138, 16).Value = 37
$ws.Cells.Item(138, 17).Value = 100
$ws.Cells.Item(138, 18).Value = "Hortaliza"
